$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add header for new column J
$ws.Range("J1").Value = "n"

# Values for column J (rows 2-13), pattern repeats 5,4,3
$values = @(5,4,3,5,4,3,5,4,3,5,4,3)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 10).Value = $values[$i]
}

# Update the selected cell/range to C4 as in the diff
$ws.Range("C4").Select()
